# ScrumLogAppDev04.xlsx edit
# "Wat ik ga doen geremoved... cause kevin couldn't handle it"
#
# - Remove the leading "wa ik ga doen: " from the D4 cell so it just reads
#   "helpen met use case diagram + wireframes opstellen"
# - Move the active selection from I4 to F4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the "wa ik ga doen: " prefix from D4's text
$ws.Range("D4").Value = "helpen met use case diagram + wireframes opstellen"

# Update the sheet's active selection/cell
$ws.Range("F4").Select()
